$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume snapshot values.
# Values that look like plain numbers (e.g. "1.000", "288.03") are written with a
# leading text-qualifying apostrophe so Excel keeps them as literal text (matching the
# original inline-string cells) instead of silently parsing them into numbers; the cell
# style is then reset to 'Normal' so no visible formatting change is introduced.
$updates = @(
    @('D2', '22.474.69'),
    @('E2', '  +0.39%  '),
    @('D3', '1.575.45'),
    @('E3', '  +0.95%  '),
    @('E4', '  -0.09%  '),
    @('E5', '  -0.14%  '),
    @('D6', '288.03'),
    @('E6', '  +0.66%  '),
    @('D7', '0.3699'),
    @('E7', '  +1.49%  '),
    @('D8', '47.72'),
    @('E8', '  -1.95%  '),
    @('E9', '  -0.13%  '),
    @('E10', '  +2.47%  '),
    @('D11', '0.07573'),
    @('E12', '  -0.11%  '),
    @('D13', '20.85'),
    @('E13', '  +0.65%  '),
    @('D14', '5.952'),
    @('E14', '  +0.86%  '),
    @('D15', '6.951'),
    @('E15', '  +1.58%  '),
    @('D16', '1.570.38'),
    @('E16', '  +0.63%  '),
    @('D17', '0.00001122'),
    @('D18', '88.33'),
    @('E18', '  -0.34%  '),
    @('D19', '0.06734'),
    @('E19', '  +0.18%  '),
    @('D20', '1.000'),
    @('D21', '6.393'),
    @('E21', '  +1.47%  '),
    @('D22', '16.54'),
    @('E22', '  +3.62%  '),
    @('E23', '  +1.23%  '),
    @('D24', '22.472.02'),
    @('E24', '  +0.39%  '),
    @('D25', '2.387'),
    @('E25', '  -0.02%  '),
    @('D26', '2.643'),
    @('E26', '  +3.24%  '),
    @('D27', '151.12'),
    @('E27', '  +1.31%  '),
    @('D28', '19.69'),
    @('E28', '  +1.72%  '),
    @('D29', '4.996'),
    @('E29', '  -0.34%  '),
    @('D30', '125.64'),
    @('E30', '  +2.40%  '),
    @('D31', '1.750.77'),
    @('E31', '  +0.78%  '),
    @('D32', '1.098'),
    @('E32', '  +4.31%  '),
    @('D33', '6.111'),
    @('E33', '  +0.29%  '),
    @('D34', '1.989'),
    @('E34', '  +0.01%  '),
    @('D35', '9.881'),
    @('E35', '  +3.47%  '),
    @('D36', '0.08367'),
    @('E36', '  +1.82%  '),
    @('E37', '  +4.18%  '),
    @('D38', '0.2242'),
    @('E38', '  +1.50%  '),
    @('E39', '  +0.89%  '),
    @('D40', '1.295'),
    @('E40', '  +0.19%  '),
    @('D41', '5.362'),
    @('E41', '  +1.16%  '),
    @('D42', '11.49'),
    @('E42', '  +3.64%  '),
    @('E43', '  +4.08%  '),
    @('D44', '14.13'),
    @('E44', '  +4.01%  '),
    @('E45', '  -0.12%  '),
    @('D46', '0.6116'),
    @('E46', '  +7.17%  '),
    @('D47', '3.779'),
    @('E47', '  +0.60%  '),
    @('D49', '125.27'),
    @('E49', '  +0.60%  '),
    @('E50', '  +0.55%  '),
    @('D51', '0.07224')
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    $bare = $newVal.Trim()
    if ($bare.Length -gt 0 -and ($bare[0] -eq "+" -or $bare[0] -eq "-")) {
        $bare = $bare.Substring(1)
    }
    # A "plain number" string (optionally signed, at most one decimal point, e.g. "1.000",
    # "288.03") would otherwise be silently parsed into a numeric cell value by Excel,
    # which both loses the original text formatting (trailing zeros, leading zeros on tiny
    # values like "0.00001122") and flips the cell's stored type away from text. Guard those
    # with a text-qualifying leading apostrophe, then restore the default "Normal" style so
    # no stray formatting is left behind on the cell.
    $isPlainNumber = ($bare.Length -gt 0) -and ($bare -match "^[0-9]+(\.[0-9]+)?$")
    if ($isPlainNumber) {
        $ws.Range($cellRef).Value = "'" + $newVal
        $ws.Range($cellRef).Style = 'Normal'
    } else {
        $ws.Range($cellRef).Value = $newVal
    }
}
